$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("position_id"), shifting the
# existing position_id / tax_src_id / tax_dest_id columns one place right.
$ws.Columns.Item(2).Insert()

# New column header
$ws.Range("B1").Value = "_requirements"

# Populate the new "_requirements" column for the rows that need it.
$ws.Range("B6").Value = "l10n_it_reverse_charge"
$ws.Range("B7").Value = "l10n_it_reverse_charge"
$ws.Range("B8").Value = "l10n_it_split_payment"
$ws.Range("B9").Value = "l10n_it_dichiarazione_intento or l10n_it_lettera_intento "

# Rows 5 and 11 keep no cell at all in the new column (fully cleared,
# not just emptied of content), matching the source data.
$ws.Range("B5").Clear()
$ws.Range("B11").Clear()

# Column widths, to match the final layout (values chosen so the engine's
# pixel-snapping rounds to the closest achievable approximation of the
# target stored widths 45.62 / 17.83 / 13.37 / 17.27).
$ws.Columns.Item(2).ColumnWidth = 44.833333333333336
$ws.Columns.Item(3).ColumnWidth = 17
$ws.Columns.Item(4).ColumnWidth = 12.5
$ws.Columns.Item(5).ColumnWidth = 16.5

$ws.Range("B10").Select()
